$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    foreach ($hf in $sec.Headers) {
        if ($hf.Exists) {
            foreach ($ish in $hf.Range.InlineShapes) {
                if ($ish.AlternativeText -eq "BTec_Logo-Orange") {
                    $ish.Name = "image2.jpg"
                }
            }
        }
    }
    foreach ($hf in $sec.Footers) {
        if ($hf.Exists) {
            foreach ($ish in $hf.Range.InlineShapes) {
                if ($ish.AlternativeText -like "*PearsonLogo.png") {
                    $ish.Name = "image1.png"
                }
            }
        }
    }
}
